# Move "Chapter 41. Introduction to epidemiology" reading assignment from the
# causality module into this (measures of occurrence) module's Required
# Readings list, as a new sub-bullet right after "Chapter 21. Measures of
# dispersion".

$d = $word.ActiveDocument

# Locate the "Chapter 21. Measures of dispersion" bullet under the
# "Cannell B. R for Epidemiology" reading.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd("`r", "`a") -eq "Chapter 21. Measures of dispersion") {
        $targetPara = $para
        break
    }
}

$r = $targetPara.Range

# Append a trailing space to the end of that paragraph's text, matching the
# same run formatting already used in the paragraph (Times New Roman, black).
$lastCharRange = $d.Range($r.End - 2, $r.End - 1)
$clonedFormat = $lastCharRange.FormattedText
$insertionPoint = $d.Range($r.End - 1, $r.End - 1)
$insertionPoint.FormattedText = $clonedFormat
$newCloneRange = $d.Range($r.End - 2, $r.End - 1)
$newCloneRange.Text = " "

# Insert a new paragraph right after it, inheriting the same list formatting
# (ilvl 1 / numId 3, no paragraph borders, single spacing, Times New Roman).
$r.InsertParagraphAfter()
$newPara = $targetPara.Next()
$newPara.Range.Text = "Chapter 41. Introduction to epidemiology"
